# Regenerate save_data for column G ("K" = strikeouts) using updated values.
# This mirrors a regen of the underlying box-score calc (K instead of Strike#)
# which changed the K column values for each of the 15 game rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 4
    4  = 10
    5  = 9
    6  = 4
    7  = 4
    8  = 8
    9  = 1
    10 = 7
    11 = 5
    12 = 4
    13 = 3
    14 = 6
    15 = 2
    16 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
